$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KPI")
$ws.Rows.Item(2).Delete()
